$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.85
$ws.Range("H2").Value = 3.5
$ws.Range("I2").Value = 3.7
$ws.Range("J2").Value = 2.5
$ws.Range("K2").Value = 2.4
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 15
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 1.6
$ws.Range("S2").Value = 1.29
$ws.Range("T2").Value = 3.5
$ws.Range("U2").Value = 1.53
$ws.Range("V2").Value = 2.38
$ws.Range("W2").Value = 10
$ws.Range("X2").Value = 11
$ws.Range("Y2").Value = 8.5
$ws.Range("Z2").Value = 17
$ws.Range("AA2").Value = 13
$ws.Range("AB2").Value = 21
$ws.Range("AC2").Value = 15
$ws.Range("AD2").Value = 7.5
$ws.Range("AE2").Value = 12
$ws.Range("AF2").Value = 34
$ws.Range("AH2").Value = 15
$ws.Range("AI2").Value = 23
$ws.Range("AJ2").Value = 13
$ws.Range("AK2").Value = 41
$ws.Range("AL2").Value = 26
$ws.Range("AM2").Value = 29
$ws.Range("AN2").Value = 4.33
$ws.Range("AO2").Value = 9.5
$ws.Range("AP2").Value = 17
$ws.Range("AQ2").Value = 29
$ws.Range("AR2").Value = 41
$ws.Range("AS2").Value = 101
$ws.Range("AT2").Value = 3.5
$ws.Range("AU2").Value = 7.5
$ws.Range("AV2").Value = 41
$ws.Range("AW2").Value = 6
$ws.Range("AX2").Value = 19
$ws.Range("AY2").Value = 23
$ws.Range("AZ2").Value = 51
$ws.Range("BA2").Value = 67
$ws.Range("BB2").Value = 126

# Row 5
$ws.Range("O5").Value = 1.2
$ws.Range("P5").Value = 4.33
$ws.Range("Q5").Value = 1.67
$ws.Range("R5").Value = 2.15

# Row 7
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 5.3
$ws.Range("J7").Value = 2.18
$ws.Range("K7").Value = 2.12
$ws.Range("L7").Value = 5.3
$ws.Range("S7").Value = 1.39
$ws.Range("T7").Value = 2.55
$ws.Range("U7").Value = 1.91
$ws.Range("X7").Value = 7.1
$ws.Range("Z7").Value = 12.5
$ws.Range("AB7").Value = 30
$ws.Range("AD7").Value = 6.7
$ws.Range("AF7").Value = 90
$ws.Range("AI7").Value = 32
$ws.Range("AK7").Value = 110
$ws.Range("AQ7").Value = 26
$ws.Range("AT7").Value = 2.55
$ws.Range("AU7").Value = 7.4
$ws.Range("AW7").Value = 6.8
$ws.Range("AZ7").Value = 200

# Row 8
$ws.Range("T8").Value = 3.25
